# --- edit.ps1 ---
# Rebuilds the eco-club notification upload summary after a data refresh:
# district-level "Eco-Club Notification Uploaded" counts increased for many
# districts, so the percentage column, Top 10 / Bottom 10 leaderboards, and
# the Overall Summary totals all need to be recalculated/rewritten to match.

$wb = $excel.ActiveWorkbook
$wsAll = $wb.Worksheets.Item("All Districts")
$wsTop = $wb.Worksheets.Item("Top 10")
$wsBottom = $wb.Worksheets.Item("Bottom 10")
$wsSummary = $wb.Worksheets.Item("Overall Summary")

# ---------------------------------------------------------------------
# 1) "All Districts": apply the updated "Eco-Club Notification Uploaded"
#    counts, then recompute "Percentage (%)" = ROUND(Uploaded/Total*100,2)
#    for every data row (rows 2-77, row 77 is the TOTAL row).
# ---------------------------------------------------------------------
$newUploadedByRow = @{
    2 = 180
    4 = 243
    5 = 128
    6 = 92
    7 = 98
    9 = 98
    10 = 157
    13 = 141
    14 = 220
    15 = 132
    16 = 123
    18 = 83
    20 = 192
    26 = 111
    27 = 251
    28 = 243
    29 = 153
    32 = 114
    33 = 296
    35 = 122
    36 = 167
    37 = 89
    39 = 299
    40 = 191
    41 = 188
    42 = 126
    44 = 320
    46 = 91
    47 = 255
    50 = 301
    51 = 179
    52 = 78
    55 = 61
    56 = 286
    57 = 143
    58 = 205
    59 = 180
    60 = 136
    61 = 236
    64 = 172
    65 = 265
    67 = 62
    69 = 120
    72 = 166
    73 = 82
    74 = 222
    75 = 378
    76 = 228
    77 = 11687
}

for ($row = 2; $row -le 77; $row++) {
    if ($newUploadedByRow.ContainsKey($row)) {
        $wsAll.Cells.Item($row, 3).Value2 = $newUploadedByRow[$row]
    }
    $totalSchools = $wsAll.Cells.Item($row, 2).Value2
    $uploaded = $wsAll.Cells.Item($row, 3).Value2
    $wsAll.Cells.Item($row, 4).Value2 = [Math]::Round(($uploaded / $totalSchools) * 100, 2)
}

# ---------------------------------------------------------------------
# 2) "Top 10": the 10 districts with the highest percentage, descending.
# ---------------------------------------------------------------------
$topRows = @(
    @("UNNAO", 435, 378),
    @("SHAMLI (PRABUDH NAGAR)", 150, 120),
    @("HAPUR (PANCHSHEEL NAGAR)", 163, 122),
    @("KHERI", 365, 255),
    @("MAHOBA", 112, 78),
    @("CHANDAULI", 290, 192),
    @("SHRAWASTI", 130, 85),
    @("RAMPUR", 266, 172),
    @("JHANSI", 300, 191),
    @("CHITRAKOOT", 134, 85),
)

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $row = $i + 2
    $item = $topRows[$i]
    $district = $item[0]
    $total = $item[1]
    $uploaded = $item[2]
    $wsTop.Cells.Item($row, 1).Value2 = $district
    $wsTop.Cells.Item($row, 2).Value2 = $total
    $wsTop.Cells.Item($row, 3).Value2 = $uploaded
    $wsTop.Cells.Item($row, 4).Value2 = [Math]::Round(($uploaded / $total) * 100, 2)
}

# ---------------------------------------------------------------------
# 3) "Bottom 10": the 10 districts with the lowest percentage, ascending.
# ---------------------------------------------------------------------
$bottomRows = @(
    @("KANPUR DEHAT", 335, 28),
    @("GHAZIPUR", 949, 91),
    @("AZAMGARH", 838, 98),
    @("MAU", 522, 61),
    @("PRAYAGRAJ", 1181, 183),
    @("AGRA", 1077, 180),
    @("KUSHINAGAR", 413, 71),
    @("MAINPURI", 516, 89),
    @("DEORIA", 580, 112),
    @("SANT KABIR NAGAR", 310, 62),
)

for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $row = $i + 2
    $item = $bottomRows[$i]
    $district = $item[0]
    $total = $item[1]
    $uploaded = $item[2]
    $wsBottom.Cells.Item($row, 1).Value2 = $district
    $wsBottom.Cells.Item($row, 2).Value2 = $total
    $wsBottom.Cells.Item($row, 3).Value2 = $uploaded
    $wsBottom.Cells.Item($row, 4).Value2 = [Math]::Round(($uploaded / $total) * 100, 2)
}

# ---------------------------------------------------------------------
# 4) "Overall Summary": total notifications uploaded + overall percentage.
# ---------------------------------------------------------------------
$grandTotalSchools = $wsAll.Cells.Item(77, 2).Value2
$grandTotalUploaded = $wsAll.Cells.Item(77, 3).Value2
$wsSummary.Cells.Item(3, 2).Value2 = $grandTotalUploaded
$wsSummary.Cells.Item(4, 2).Value2 = [Math]::Round(($grandTotalUploaded / $grandTotalSchools) * 100, 2)

